$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3852.8667
$ws.Range("I19").Value = 3660.7368
$ws.Range("J19").Value = 4184.727
$ws.Range("K19").Value = 3660.7368
$ws.Range("L19").Value = 4184.727
$ws.Range("M19").Value = -3485.7368
$ws.Range("N19").Value = -4534.727

$ws.Range("H62").Value = 6810.0557
$ws.Range("I62").Value = 5869.4546
$ws.Range("K62").Value = 5869.4546
$ws.Range("M62").Value = -5245.4546

$ws.Range("H65").Value = 6810.0557
$ws.Range("I65").Value = 5869.4546
$ws.Range("K65").Value = 29347.273
$ws.Range("M65").Value = -26227.273

$ws.Range("H106").Value = 1658.1666
$ws.Range("I106").Value = 1544.2222
$ws.Range("K106").Value = 1544.2222
$ws.Range("M106").Value = -913.2221999999999

$ws.Range("H116").Value = 5975.619
$ws.Range("J116").Value = 6089.846
$ws.Range("L116").Value = 6089.846
$ws.Range("N116").Value = -12973.846

$ws.Range("H125").Value = 32017.2
$ws.Range("J125").Value = 2825
$ws.Range("L125").Value = 25425
$ws.Range("N125").Value = -30345

$ws.Range("H137").Value = 20155.234
$ws.Range("I137").Value = 26899.166
$ws.Range("K137").Value = 80697.498
$ws.Range("M137").Value = -78147.498

$ws.Range("H141").Value = 1695.4286
$ws.Range("I141").Value = 1575.8
$ws.Range("J141").Value = 1994.5
$ws.Range("K141").Value = 4727.4
$ws.Range("L141").Value = 5983.5
$ws.Range("M141").Value = 452.6000000000004
$ws.Range("N141").Value = -16343.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = $null

$ws.Range("H45").Value = 3213.9167
$ws.Range("I45").Value = 1754.5454
$ws.Range("K45").Value = 1754.5454
$ws.Range("M45").Value = -1377.5454

$ws.Range("H61").Value = 1695.2
$ws.Range("I61").Value = 1226.4706
$ws.Range("K61").Value = 1226.4706
$ws.Range("M61").Value = -1014.4706

$ws.Range("H74").Value = 546629.8
$ws.Range("I74").Value = 667337
$ws.Range("K74").Value = 667337
$ws.Range("M74").Value = -666463

$ws.Range("H77").Value = 546629.8
$ws.Range("I77").Value = 667337
$ws.Range("K77").Value = 3336685
$ws.Range("M77").Value = -3332317

$ws.Range("H97").Value = 1422.5
$ws.Range("I97").Value = 1009.64
$ws.Range("K97").Value = 1009.64
$ws.Range("M97").Value = -513.64

$ws.Range("H136").Value = 1695.2
$ws.Range("I136").Value = 1226.4706
$ws.Range("K136").Value = 3679.4118
$ws.Range("M136").Value = -1129.4118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 110780
$ws.Range("J59").Value = 110780
$ws.Range("L59").Value = 110780
$ws.Range("N59").Value = -112474

$ws.Range("H107").Value = 19825
$ws.Range("I107").Value = 23319.088
$ws.Range("K107").Value = 23319.088
$ws.Range("M107").Value = -21399.088

$ws.Range("H134").Value = 3142.1052
$ws.Range("I134").Value = 3038.889
$ws.Range("K134").Value = 9116.667000000001
$ws.Range("M134").Value = -6581.667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3847979.8
$ws.Range("I31").Value = 4546971.5
$ws.Range("J31").Value = 3525.75
$ws.Range("K31").Value = 4546971.5
$ws.Range("L31").Value = 3525.75
$ws.Range("M31").Value = -4546676.5
$ws.Range("N31").Value = -4115.75

$ws.Range("H34").Value = 3847979.8
$ws.Range("I34").Value = 4546971.5
$ws.Range("J34").Value = 3525.75
$ws.Range("K34").Value = 4546971.5
$ws.Range("L34").Value = 3525.75
$ws.Range("M34").Value = -4546769.5
$ws.Range("N34").Value = -3929.75

$ws.Range("H94").Value = 1654
$ws.Range("J94").Value = 1070
$ws.Range("L94").Value = 1070
$ws.Range("N94").Value = -1972

$ws.Range("H105").Value = 2175.7334
$ws.Range("I105").Value = 1148.7273
$ws.Range("K105").Value = 1148.7273
$ws.Range("M105").Value = 598.2727

$ws.Range("H132").Value = 25654.12
$ws.Range("I132").Value = 31304
$ws.Range("K132").Value = 93912
$ws.Range("M132").Value = -91382

$ws.Range("H134").Value = 2385.6667
$ws.Range("I134").Value = 1957.1765
$ws.Range("J134").Value = 4206.75
$ws.Range("K134").Value = 5871.529500000001
$ws.Range("L134").Value = 12620.25
$ws.Range("M134").Value = -3336.529500000001
$ws.Range("N134").Value = -17690.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1080.25
$ws.Range("I5").Value = 1091.7142
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 3275.1426
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -3163.1426
$ws.Range("N5").Value = -3224

$ws.Range("H113").Value = 3802.1667
$ws.Range("J113").Value = 3802.1667
$ws.Range("L113").Value = 11406.5001
$ws.Range("N113").Value = -15746.5001

$ws.Range("H131").Value = 329841.7
$ws.Range("I131").Value = 1418074
$ws.Range("J131").Value = 3372
$ws.Range("K131").Value = 4254222
$ws.Range("L131").Value = 10116
$ws.Range("M131").Value = -4249182
$ws.Range("N131").Value = -20196

$ws.Range("H132").Value = 1497.1818
$ws.Range("I132").Value = 1528.5
$ws.Range("J132").Value = 1459.6
$ws.Range("K132").Value = 13756.5
$ws.Range("L132").Value = 13136.4
$ws.Range("M132").Value = -11226.5
$ws.Range("N132").Value = -18196.4

$ws.Range("H135").Value = 1080.25
$ws.Range("I135").Value = 1091.7142
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9825.4278
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -7290.427799999999
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 51840.668
$ws.Range("J96").Value = 51840.668
$ws.Range("L96").Value = 51840.668
$ws.Range("N96").Value = -57332.668

$ws.Range("H132").Value = 3506
$ws.Range("I132").Value = 3506
$ws.Range("K132").Value = 10518
$ws.Range("M132").Value = -7988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 10727.6
$ws.Range("J104").Value = 10727.6
$ws.Range("L104").Value = 10727.6
$ws.Range("N104").Value = -17715.6

$ws.Range("H132").Value = 3945.8
$ws.Range("I132").Value = 3945.8
$ws.Range("K132").Value = 11837.4
$ws.Range("M132").Value = -9307.400000000001

$ws.Range("H136").Value = 5254.035
$ws.Range("I136").Value = 4360.256
$ws.Range("J136").Value = 7999.2144
$ws.Range("K136").Value = 13080.768
$ws.Range("L136").Value = 23997.6432
$ws.Range("M136").Value = -10530.768
$ws.Range("N136").Value = -29097.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4195.5713
$ws.Range("I132").Value = 4268.7896
$ws.Range("K132").Value = 12806.3688
$ws.Range("M132").Value = -10276.3688

$ws.Range("H136").Value = 12289.981
$ws.Range("I136").Value = 13132.412
$ws.Range("J136").Value = 1549
$ws.Range("K136").Value = 39397.236
$ws.Range("L136").Value = 4647
$ws.Range("M136").Value = -36847.236
$ws.Range("N136").Value = -9747
